# Update "想去人数" (want-to-go count) figures in column F for the
# sheets "展览" and "全部类型" to reflect the latest scrape results.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 2955
    "F5"  = 6718
    "F6"  = 1693
    "F10" = 117
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
